# Applies the "Dang ky, tiep nhan, kiem tra du lieu dau vao, tai file" edit.
# Uses Find/Replace on unique, context-rich phrases so each substitution
# lands on the correct occurrence even though several short numeric tokens
# (08, 15, 19, 23 ...) and names (Ma Kien Tu / Ha The Duy) repeat in the
# document.

$d = $word.ActiveDocument

function Replace-Unique($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $findText"
    }
    return $ok
}

# 1) "Hoi 15 gio 08 phut, ngay 19 thang ..." -> "Hoi 09 gio 27 phut, ngay 23 thang ..."
Replace-Unique "Hồi 15 giờ 08 phút, ngày 19 tháng" "Hồi 09 giờ 27 phút, ngày 23 tháng"

# 2) Quyet dinh trung cau giam dinh so
Replace-Unique "185/QĐ-CQĐT" "193/QĐ-TCGĐ"

# 3) ngay cua quyet dinh
Replace-Unique "18/05/2023 " "23/05/2023 "

# 4) Ben giao: Ha The Duy
Replace-Unique "Bên giao: Hà Thế Duy" "Bên giao: Chu An Khánh"

# 5) Ben nhan: Ma Kien Tu;
Replace-Unique "Bên nhận: Ma Kiên Tú;" "Bên nhận: Trần Xuân Quang;"

# 6) Chuc vu cua ben nhan: "Cán bộ" (no trailing space) -> "Giám định viên"
#    NB: the "Bên giao" role text is "Cán bộ " (with a trailing space) and
#    must stay untouched, so search only the part of the document after
#    "Bên nhận" to uniquely hit the second occurrence.
$anchor = $d.Content
$null = $anchor.Find.Execute("Bên nhận", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterAnchor = $anchor.End
$sub = $d.Range($afterAnchor, $d.Content.End)
$ok = $sub.Find.Execute("Cán bộ", $true, $false, $false, $false, $false, $true, 1, $false, "Giám định viên", 2)
if (-not $ok) { Write-Host "NOT FOUND: Cán bộ (Bên nhận)" }

# 7) Doi tuong giam dinh mo ta (the whole paragraph's text is replaced)
$oldObj = "- 01 đoạn video có tên " + [char]0x201C + "ch0_20230416101915_013.mp4" + [char]0x201D + ", mã MD5: 8bc02026d02c347db852c64d425468d2 được lưu trữ trong thẻ nhớ có chữ Pioneer 32GB, niêm phong trong 01 bì thư ghi " + [char]0x22 + "Thẻ nhớ camera hành trình xe ô tô BKS: 000.46"
$newObj = "01 (một) điện thoại di động nhãn hiệu Samsung Galaxy Note 9, màu xanh dương đen, có số Imei: 352141101123160, Imei 2: 352142101123168, điện thoại đã qua sử dụng, có sọc kẻ màu xanh trên màn hình."
Replace-Unique $oldObj $newObj

# 8) "Viec giao, nhan ket thuc hoi 15 gio 23 phut, ngay 19 thang ..." -> "... hoi 09 gio 42 phut, ngay 23 thang ..."
Replace-Unique "Việc giao, nhận kết thúc hồi 15 giờ 23 phút, ngày 19 tháng" "Việc giao, nhận kết thúc hồi 09 giờ 42 phút, ngày 23 tháng"

# 9) Signature table - "BEN NHAN" signatory name
Replace-Unique "Ma Kiên Tú" "Trần Xuân Quang"

# 10) Signature table - "BEN GIAO" signatory name
Replace-Unique "Hà Thế Duy" "Chu An Khánh"

Write-Host "Done"
